$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper pattern used throughout this script:
#   The engine merges adjacent runs that share identical formatting
#   whenever a Range.Text assignment/insert touches the document, so a
#   naive Find&Replace would fuse the edited text back into its
#   neighbouring (unrelated) runs. To keep the run layout the diff
#   expects, we first drop temporary "wall" bookmarks at every offset
#   where a run boundary must survive, confined to a pre-existing
#   NON-EMPTY span of original text (re-writing sub-ranges of existing
#   text keeps the walls in place, unlike inserting into a zero-width
#   gap, which the engine always merges into the run on its left).
#   Once every chunk has been rewritten we delete the temporary walls
#   (any permanent bookmark required by the diff is added separately).
# ------------------------------------------------------------------

# ====================================================================
# Edit 1: <meta name="title" ...> text
#   "13 Door System" + ' ">'  ->  "15 Animation for Doors" + '">'
# ====================================================================

$rngA = $d.Content
$rngA.Find.Execute('13 Door System') | Out-Null
$aStart = $rngA.Start
$aEnd = $rngA.End

$rngB = $d.Content
$rngB.Find.Execute(' "' + [char]62) | Out-Null
$bEnd = $rngB.End

# Wall bookmarks at every boundary that must survive as a run edge:
# before chunk A, between A and B, and after chunk B. Placed up front
# (while all the original text is still intact) so later rewrites of
# existing sub-ranges cannot bridge across them.
$d.Bookmarks.Add('ZZWALL1A', $d.Range($aStart, $aStart)) | Out-Null
$d.Bookmarks.Add('ZZWALL1M', $d.Range($aEnd, $aEnd)) | Out-Null
$d.Bookmarks.Add('ZZWALL1B', $d.Range($bEnd, $bEnd)) | Out-Null

# Rewrite chunk A (between walls 1A and 1M); re-read the live bookmark
# offsets rather than trusting the pre-edit numbers, since earlier
# rewrites shift everything after them.
$w1a = $d.Bookmarks('ZZWALL1A').Start
$w1m = $d.Bookmarks('ZZWALL1M').Start
$d.Range($w1a, $w1m).Text = '15 Animation for Doors'

# Rewrite chunk B (between walls 1M and 1B) using freshly re-read offsets.
$w1m = $d.Bookmarks('ZZWALL1M').Start
$w1b = $d.Bookmarks('ZZWALL1B').Start
$d.Range($w1m, $w1b).Text = '"' + [char]62

$d.Bookmarks('ZZWALL1A').Delete()
$d.Bookmarks('ZZWALL1M').Delete()
$d.Bookmarks('ZZWALL1B').Delete()

# ====================================================================
# Edit 2: <meta name="description" ...> text
#   Split the old sentence into two runs wrapped by a permanent
#   bookmark (_Hlk219378004).
# ====================================================================

$rngFull = $d.Content
$rngFull.Find.Execute('In this tutorial, we will be starting to build our door system. In the Dungeon Crawler game, we will want to have a variety of different doors, which can be coded to take the hero to various areas, in the game.') | Out-Null
$fullStart = $rngFull.Start
$fullEnd = $rngFull.End

$midPos = $fullStart + 18   # length of 'In this tutorial, '

$d.Bookmarks.Add('ZZWALL2A', $d.Range($fullStart, $fullStart)) | Out-Null
$d.Bookmarks.Add('ZZWALL2M', $d.Range($midPos, $midPos)) | Out-Null
$d.Bookmarks.Add('ZZWALL2B', $d.Range($fullEnd, $fullEnd)) | Out-Null

$w2a = $d.Bookmarks('ZZWALL2A').Start
$w2m = $d.Bookmarks('ZZWALL2M').Start
$d.Range($w2a, $w2m).Text = 'In this tutorial, '

$w2m = $d.Bookmarks('ZZWALL2M').Start
$w2b = $d.Bookmarks('ZZWALL2B').Start
$d.Range($w2m, $w2b).Text = 'we will be looking at writing the code to give our doors some special effects. Yes, we will be animating those effects right in Game Maker' + [char]8217 + 's code panel and giving them a bit of razzle-dazzle.'

$bm2Start = $d.Bookmarks('ZZWALL2A').Start
$bm2End = $d.Bookmarks('ZZWALL2B').Start
$d.Bookmarks.Add('_Hlk219378004', $d.Range($bm2Start, $bm2End)) | Out-Null

$d.Bookmarks('ZZWALL2A').Delete()
$d.Bookmarks('ZZWALL2M').Delete()
$d.Bookmarks('ZZWALL2B').Delete()

# ====================================================================
# Edit 3: <meta name="url" ...> text
#   Insert a standalone run holding a single space before the filename
#   run, and update the filename run's path segment.
# ====================================================================

$rngUrl = $d.Content
$rngUrl.Find.Execute('Enlightenment/Articles/2026/2_Game_Maker_2/13_Door_System/13_Door_System.html') | Out-Null
$urlStart = $rngUrl.Start
$urlEnd = $rngUrl.End

$urlMidPos = $urlStart + 1   # first character becomes the new space-only run

$d.Bookmarks.Add('ZZWALL3A', $d.Range($urlStart, $urlStart)) | Out-Null
$d.Bookmarks.Add('ZZWALL3M', $d.Range($urlMidPos, $urlMidPos)) | Out-Null
$d.Bookmarks.Add('ZZWALL3B', $d.Range($urlEnd, $urlEnd)) | Out-Null

$w3a = $d.Bookmarks('ZZWALL3A').Start
$w3m = $d.Bookmarks('ZZWALL3M').Start
$d.Range($w3a, $w3m).Text = ' '

$w3m = $d.Bookmarks('ZZWALL3M').Start
$w3b = $d.Bookmarks('ZZWALL3B').Start
$d.Range($w3m, $w3b).Text = 'Enlightenment/Articles/2026/2_Game_Maker_2/15_Animation_For_Doors/15_Animation_For_Doors.html'

$d.Bookmarks('ZZWALL3A').Delete()
$d.Bookmarks('ZZWALL3M').Delete()
$d.Bookmarks('ZZWALL3B').Delete()

Write-Output 'done'
